$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was the generic "1", now named after the municipality)
$ws.Name = "ქობულეთი"

# Drop the obsolete "(census results)" sub-heading row — everything below
# it shifts up by one row.
$ws.Rows("2").Delete()

# The table used to show three census years (1989 / 2002 / 2014) side by
# side; keep only the most recent (2014) figures and drop the other two
# year columns.
$ws.Columns("B:C").Delete()

# Restore the active-cell selection left by the author when they saved.
$ws.Range("A2").Select() | Out-Null
